# "Improved Query error reporting."
#
# 1) word/document.xml : remove the three demo paragraphs
#    ("A simple demonstration of a query :", the "m:self." field-code
#    paragraph, and "End of demonstration.") that used to precede the
#    bookmark paragraph, leaving only the (now first) bookmark paragraph.
#
# 2) word/header1.xml : collapse the first header paragraph's five runs
#    (split apart by spell-check <w:proofErr> markers) into two runs -
#    "A simple demonstration of a " and "query\u00a0:" - tag the language
#    as en-US on the paragraph mark and on both runs, and replace the
#    spell-check proofErr pair with a single gramStart/gramEnd pair
#    wrapping "query\u00a0:".

$d = $word.ActiveDocument

# --- 1) document.xml --------------------------------------------------
# Paragraph 3 is "End of demonstration."; deleting from the very start of
# the document through the end of that paragraph (its mark included)
# removes the three leading paragraphs and leaves the bookmark paragraph
# as the new first (and only) paragraph.
$lastDemoPara = $d.Paragraphs.Item(3)
$removeRange = $d.Range(0, $lastDemoPara.Range.End)
$removeRange.Delete()

# --- 2) header1.xml -----------------------------------------------------
$hdr = $d.Sections.Item(1).Headers.Item(1)

# The first header paragraph currently reads (5 runs, nbsp before the
# colon): "A simple " + "demonstration" + " of a " + "query" + "\u00a0:"
$firstPara = $hdr.Range.Duplicate
$firstPara.Collapse(1)
[void]$firstPara.MoveEnd(1, 35)

$nbsp = [char]0xA0
$newParaXml = @"
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">A simple demonstration of a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>query${nbsp}:</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$firstPara.InsertXML($newParaXml)
